$d = $word.ActiveDocument

$replacements = @(
    @("304×8=", "336×3="),
    @("450×3=", "570×6="),
    @("876×8=", "674×3="),
    @("333×9=", "854×5="),
    @("323×6=", "278×9="),
    @("845×8=", "954×7="),
    @("947×3=", "327×2="),
    @("833×4=", "564×6="),
    @("951×5=", "123×3="),
    @("572×3=", "442×6="),
    @("341×6=", "165×8="),
    @("787×6=", "785×5="),
    @("389×4=", "989×3="),
    @("483×8=", "361×8="),
    @("742×4=", "635×9="),
    @("453×7=", "831×8="),
    @("547×2=", "264×7="),
    @("907×7=", "678×9="),
    @("350×9=", "301×6="),
    @("837×3=", "963×2="),
    @("493×5=", "246×4="),
    @("640×3=", "419×7="),
    @("891×7=", "553×8="),
    @("961×7=", "408×3="),
    @("569×3=", "680×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
